$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) values are stored as literal text (e.g. "68.018.13", "10.60")
# in the source data, not as real numbers. Force text storage via NumberFormat
# "@" before assignment so Excel does not reinterpret/round them as doubles,
# then restore the default "Normal" style so no stray formatting is left behind.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.049.79"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.10%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.784.38"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.82%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.81"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.75%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "170.07"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.19%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.782.71"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.68%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("E9").Value = "  -0.31%  "

$ws.Range("E10").Value = "  -2.09%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.53"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.02%  "

$ws.Range("E12").Value = "  -1.56%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000279"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.03%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.75"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.55%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.414.16"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.72%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.789.30"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.97%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.79"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.72%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "67.955.63"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.16%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.21"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.41%  "

$ws.Range("E20").Value = "  -0.26%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.60"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.17%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "468.24"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.19%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.721"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.06%  "

$ws.Range("E24").Value = "  -8.86%  "

$ws.Range("E25").Value = "  +0.02%  "

$ws.Range("E26").Value = "  +0.30%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.17"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.01%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.49"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.02%  "

$ws.Range("E29").Value = "  +0.01%  "

$ws.Range("E30").Value = "  -1.27%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.929.27"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.69%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.62"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.73%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "30.58"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.87%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.24"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.99%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.27"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.93%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.737.40"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.05%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.104"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.63%  "

$ws.Range("E38").Value = "  -6.67%  "

$ws.Range("E39").Value = "  -0.01%  "

$ws.Range("E40").Value = "  -1.64%  "

$ws.Range("E41").Value = "  -1.57%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.999"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.09%  "

$ws.Range("E43").Value = "  -1.08%  "

$ws.Range("E44").Value = "  +0.01%  "

$ws.Range("E45").Value = "  -0.04%  "

$ws.Range("E46").Value = "  -2.68%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "402.93"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.51%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "45.65"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.76%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.000278"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -7.50%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "40.37"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.54%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "142.04"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.32%  "

Write-Output "Applied cryptos update"